# LOLER-ReportTemplate.xlsx edit
# Commit message: "template now works for LOLER Documents"
#
# The sheet is repurposed from a Lee-Lifting-specific report (with a
# hard-coded client address block and a 14-column defect/inspection table)
# into a generic template that uses {placeholder} tokens for the client
# details and drops the now-unused summary table entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Table1" summary/defects table (rows 14:19, cols A:N) ---
# Removing the ListObject drops the table part + its autofilter/column
# definitions; deleting the now-empty rows collapses the sheet's used
# range back down to A2:D12.
if ($ws.ListObjects.Count -gt 0) {
    $ws.ListObjects.Item(1).Delete()
}
$ws.Rows("14:19").Delete()

# --- Replace the hard-coded "Lee Lifting Services" address block with
#     generic client placeholders ---
$ws.Range("A7").Value = "{Client Name}"
$ws.Range("A8").Value = "{Client Address}"
$ws.Range("A9").Value = ""
$ws.Range("A10").Value = ""
$ws.Range("A11").Value = ""
$ws.Range("A12").Value = ""

# --- Job Information panel: fold the job number into one templated
#     label, replace "Reference:" with "Contact:", and clear the other
#     now-unused labels (Job Type / Logged / Tel) ---
$ws.Range("C10").Value = "Job Number: {Job Number}"
$ws.Range("D10").Value = "Contact: "
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""

# --- Restore the selection to match the saved view state ---
$ws.Range("A29").Select()
